$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert "Icon" and "ShowName" columns before the
# existing "Desc" column, shifting Desc from G1 to I1 ---
$ws.Range("G1").Value = "Icon"
$ws.Range("H1").Value = "ShowName"
$ws.Range("I1").Value = "Desc"

# --- Data rows 2-10 ---
# Column D already holds "Prefabs/Object/<Name>_1_1"; column G (Icon) gets
# the short name (same text minus the "Prefabs/Object/" prefix), column H
# (ShowName) and column I (Desc, shifted from the old G) both get the
# Chinese description that used to live in column G.
$names = @("Altar_1_1","Arena_1_1","Camp_1_1","GoldMine_1_1","Item_hourse_1_1","League_1_1","MagicHourse_1_1","Tower_1_1","Town_1_1")
$descs = @("一级祭坛","一级竞技场","一级兵营","一级金矿","一级道具屋","一级公会","一级魔法屋","一级箭塔","一级大厅")

for ($i = 0; $i -lt 9; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 7).Value = $names[$i]
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = $descs[$i]
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 9).Value = $descs[$i]
    $ws.Cells.Item($r, 9).NumberFormat = "@"
}

# --- Column widths: new columns G:I all match the old "Desc" column width ---
$ws.Columns.Item(7).ColumnWidth = 10.29
$ws.Columns.Item(8).ColumnWidth = 10.29
$ws.Columns.Item(9).ColumnWidth = 10.29

# --- Selection moves to H10 ---
$ws.Range("H10").Select()
